# Apply "gain de pièces lors d'une victoire" edits to the Antoine.xlsx
# work-log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text updates -----------------------------------------------------
$ws.Range("B11").Value = "Discussions de groupe, travail sur wordanalyzer"
$ws.Range("B24").Value = "Endpoint powers corrigé et opérationel"

# --- Hours updates (the "Heures" column) -------------------------------
$ws.Range("C13").Value = 2
$ws.Range("C18").Value = 1.5

# C32 holds =SUM(C5:C31) and recalculates automatically to reflect
# the updated hours above.

# --- Selection / active cell -------------------------------------------
$ws.Range("B25").Select()
